# Add a new "Remove" / "Entfernen" translation row to both the "en" and
# "de" worksheets of the shared-resource workbook (row 60), matching the
# existing table's layout (language key in column A, language value in
# column B).

$wb = $excel.ActiveWorkbook

$wsEn = $wb.Worksheets.Item("en")
$wsDe = $wb.Worksheets.Item("de")

# --- "en" sheet: A60 = "Remove", B60 = "Remove" -----------------------
$wsEn.Range("A58:B58").Copy()
$wsEn.Range("A60:B60").PasteSpecial(-4122)  # xlPasteFormats
$wsEn.Application.CutCopyMode = 0

$wsEn.Range("A60").Value = "Remove"
$wsEn.Range("B60").Value = "Remove"

# --- "de" sheet: A60 = "Remove", B60 = "Entfernen" ---------------------
$wsDe.Range("A58:B58").Copy()
$wsDe.Range("A60:B60").PasteSpecial(-4122)  # xlPasteFormats
$wsDe.Application.CutCopyMode = 0

$wsDe.Range("A60").Value = "Remove"
$wsDe.Range("B60").Value = "Entfernen"

# Leave selection/cursor roughly where the original author left it after
# typing the new row.
$wsDe.Activate()
$wsDe.Range("B60").Select()

$wsEn.Activate()
$wsEn.Range("A61").Select()
